$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Tạo menu và hotkey" (row 9) finished on 22 tháng 10 - record the actual
# start/end date now that the task is done.
$ws.Range("H9").Value = "22 tháng 10"
$ws.Range("I9").Value = "22 tháng 10"

# Move the active selection to reflect where the user last worked
$ws.Range("I9").Select()
